# Generate Report for Handoff
# Updates the handoff-generation timestamps for the "d77c1228-..." row (row 5)
# across the Overview, zh-cn and de-de worksheets, as part of generating the
# handoff report.

$wb = $excel.ActiveWorkbook

# zh-cn sheet: "Latest Handoff Datetime" (column H) for row 5 gets a fresh timestamp
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H5").Value = "2016-09-05 10:22:49"

# de-de sheet: "Latest Handoff Datetime" (column H) for row 5 gets a fresh timestamp
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H5").Value = "2016-09-05 10:23:07"

# Overview sheet: "Latest HO Xliff Generate Date" (column G) for row 5 gets a fresh timestamp
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G5").Value = "2016-09-05 10:23:07"
